$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.196.15"
$ws.Range("E2").Value = "  +3.36%  "
$ws.Range("D3").Value = "2.060.22"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "230.70"
$ws.Range("E5").Value = "  +2.44%  "
$ws.Range("E6").Value = "  +1.78%  "
$ws.Range("D7").Value = "58.30"
$ws.Range("E7").Value = "  +7.23%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +3.53%  "
$ws.Range("D10").Value = "0.0809"
$ws.Range("E10").Value = "  +3.52%  "
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "2.366.74"
$ws.Range("E12").Value = "  +3.31%  "
$ws.Range("D13").Value = "14.64"
$ws.Range("E13").Value = "  +4.38%  "
$ws.Range("D14").Value = "20.68"
$ws.Range("E14").Value = "  +3.19%  "
$ws.Range("D15").Value = "0.755"
$ws.Range("E15").Value = "  +2.74%  "
$ws.Range("E16").Value = "  +4.26%  "
$ws.Range("D17").Value = "2.061.74"
$ws.Range("E17").Value = "  +3.66%  "
$ws.Range("D18").Value = "38.101.91"
$ws.Range("E18").Value = "  +3.37%  "
$ws.Range("D19").Value = "6.16"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").Value = "69.98"
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("D21").Value = "0.0₃0833"
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("D22").Value = "225.08"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "2.45"
$ws.Range("E24").Value = "  +1.48%  "
$ws.Range("E25").Value = "  +3.73%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "9.34"
$ws.Range("E26").Value = "  +3.02%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "166.81"
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("D28").Value = "0.133"
$ws.Range("E28").Value = "  +7.49%  "
$ws.Range("D29").Value = "19.07"
$ws.Range("E29").Value = "  +2.61%  "
$ws.Range("E30").Value = "  +2.44%  "
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("D32").Value = "4.57"
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("D33").Value = "4.64"
$ws.Range("E33").Value = "  +5.74%  "
$ws.Range("D34").Value = "0.0616"
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("E35").Value = "  +7.20%  "
$ws.Range("D36").Value = "2.37"
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("D37").Value = "6.12"
$ws.Range("E37").Value = "  +16.18%  "
$ws.Range("E38").Value = "  +5.99%  "
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("E40").Value = "  +2.70%  "
$ws.Range("D41").Value = "98.40"
$ws.Range("E41").Value = "  +4.33%  "
$ws.Range("D42").Value = "17.09"
$ws.Range("E42").Value = "  +4.80%  "
$ws.Range("D43").Value = "1.481.22"
$ws.Range("E43").Value = "  +1.34%  "
$ws.Range("D44").Value = "0.0943"
$ws.Range("E45").Value = "  +4.15%  "
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").Value = "4.15"
$ws.Range("E47").Value = "  +18.65%  "
$ws.Range("E48").Value = "  +1.92%  "
$ws.Range("D49").Value = "2.97"
$ws.Range("E49").Value = "  +2.42%  "
$ws.Range("D50").Value = "7.11"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "2.252.50"
$ws.Range("E51").Value = "  +3.31%  "
